$d = $word.ActiveDocument

# --- Title paragraph: "Spectralclust Package" -> "Dppsampl " + "Package" -----
# The original title run reads "Spectralclust Package". The new title keeps
# "Package" but changes the preceding word to "Dppsampl", and that trailing
# word ends up split into its own run (separate <w:r>) with the same
# character formatting as the rest of the heading.

$titlePara = $d.Paragraphs(1).Range
if ($titlePara.Text.TrimEnd([char]13, [char]7) -ne "Spectralclust Package") {
    throw "Unexpected title paragraph text: [$($titlePara.Text)]"
}

# 1) Replace the leading word (keep the trailing space) in place - this
#    keeps it as a single run with the original run formatting/rsids.
$d.Content.Find.Execute("Spectralclust ", $true, $false, $false, $false, `
    $false, $true, 1, $false, "Dppsampl ", 2) | Out-Null

# 2) Re-locate "Package" within the (now shorter) title paragraph and nudge
#    its character formatting so Word has to materialize it as a distinct
#    run, separate from the "Dppsampl " run, while keeping the same
#    rendered font (Calibri Light, for ascii/hAnsi/eastAsia/cs).
$titlePara = $d.Paragraphs(1).Range
$pkgStart = $titlePara.Start + 9
$pkgEnd = $pkgStart + 7
$pkgRange = $d.Range($pkgStart, $pkgEnd)
if ($pkgRange.Text -ne "Package") {
    throw "Unexpected range text: [$($pkgRange.Text)]"
}
$pkgRange.Font.Name = "Calibri Light"
$pkgRange.Font.NameFarEast = "Calibri Light"
$pkgRange.Font.NameOther = "Calibri Light"
$pkgRange.Font.NameBi = "Calibri Light"
